# Adds two additional response rows for "cualidades_valora_candidato"
# (it becomes a 3-answer / multi-punch question: _O1, _O2, _O3), matching
# the pattern already used elsewhere in the dictionary (e.g. problema_chile_O1/O2/O3).
#
# Before:
#   row 102: cualidades_valora_candidato
#   row 103: otro_cualidades_valora_candidato
#
# After:
#   row 102: cualidades_valora_candidato_O1
#   row 103: cualidades_valora_candidato_O2   (new)
#   row 104: cualidades_valora_candidato_O3   (new)
#   row 105: otro_cualidades_valora_candidato (shifted down by 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right after row 102, pushing everything else down.
$ws.Rows("103:104").Insert()

# Row 103 = copy of row 102 (bloque/pregunta/tipo/respuesta) with its own "llave".
$ws.Range("A103").Value = $ws.Range("A102").Value2
$ws.Range("B103").Value = $ws.Range("B102").Value2
$ws.Range("C103").Value = $ws.Range("C102").Value2
$ws.Range("D103").Value = "cualidades_valora_candidato_O2"
$ws.Range("F103").Value = $ws.Range("F102").Value2
$ws.Rows("103").AutoFit()

# Row 104 = copy of row 102 (bloque/pregunta/tipo/respuesta) with its own "llave".
$ws.Range("A104").Value = $ws.Range("A102").Value2
$ws.Range("B104").Value = $ws.Range("B102").Value2
$ws.Range("C104").Value = $ws.Range("C102").Value2
$ws.Range("D104").Value = "cualidades_valora_candidato_O3"
$ws.Range("F104").Value = $ws.Range("F102").Value2
$ws.Rows("104").AutoFit()

# Original row 102's "llave" gains the _O1 suffix.
$ws.Range("D102").Value = "cualidades_valora_candidato_O1"
